$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D ("Tipo") to make room for "MAE".
# This shifts the existing "Tipo" header and "single" value from D to E.
$ws.Columns.Item(4).Insert()

# New header for the inserted column D (copy header formatting from C1)
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "MAE"

# New MAE value for row 2
$ws.Range("D2").Value = 0.1783322799334486

# Update existing MSE (B2) and R2 (C2) values
$ws.Range("B2").Value = 0.06046199291159801
$ws.Range("C2").Value = 0.9994264909013122
